$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3:C12").ClearContents()
$ws.Range("D10").Select()
